# "Combining data quality measurement with data model quality."
#
# The workbook's 2nd..5th tabs are renamed and their tables are reshuffled
# (the underlying sheet files / relationship ids keep their identity, only
# the displayed tab name and the cell contents change):
#
#   Sheet#2 "1_ISSUES"          -> "SCHEME_MEASURES"   (now holds scheme-level measure counts)
#   Sheet#3 "2_SCHEME_MEASURES" -> "METADATA_ISSUES"   (now holds the issues table)
#   Sheet#4 "3_MODEL_MEASURES"  -> "METADATA_MEASURES" (gains two extra measure rows)
#   Sheet#5 "4_MODEL_METRICS"   -> "METADATA_METRICS"  (indicator codes/labels renumbered, one row removed)

$wb = $excel.ActiveWorkbook

$wsIssues   = $wb.Worksheets.Item(2)   # currently "1_ISSUES"
$wsScheme   = $wb.Worksheets.Item(3)   # currently "2_SCHEME_MEASURES"
$wsModelM   = $wb.Worksheets.Item(4)   # currently "3_MODEL_MEASURES"
$wsModelMet = $wb.Worksheets.Item(5)   # currently "4_MODEL_METRICS"

# ---------------------------------------------------------------------------
# Before touching any content: copy the bold/bordered header style that the
# issues sheet (D1:H1) currently has over to the sheet that will become the
# new issues table ("METADATA_ISSUES"), so its wider header keeps the same
# look as the rest of the header row.
# ---------------------------------------------------------------------------
$wsIssues.Range("D1:H1").Copy()
$wsScheme.Range("D1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet #2 : "1_ISSUES" -> "SCHEME_MEASURES"
# Replace the issues table with the scheme-level measure counts.
# ---------------------------------------------------------------------------
$wsIssues.Range("D1:H2").Clear()

$wsIssues.Range("A1").Value = "Indicator"
$wsIssues.Range("B1").Value = "Description"
$wsIssues.Range("C1").Value = "Value"

$wsIssues.Range("A2").Value = "MQMS01"
$wsIssues.Range("B2").Value = "Total number of tables"
$wsIssues.Range("C2").Value = 6

$wsIssues.Range("A3").Value = "MQMS02"
$wsIssues.Range("B3").Value = "Total number of columns"
$wsIssues.Range("C3").Value = 246

$wsIssues.Range("A4").Value = "MQMS03"
$wsIssues.Range("B4").Value = "Total number of primary key"
$wsIssues.Range("C4").Value = 3

$wsIssues.Range("A5").Value = "MQMS04"
$wsIssues.Range("B5").Value = "Total number of foreign key"
$wsIssues.Range("C5").Value = 0

$wsIssues.Range("A6").Value = "MQMS05"
$wsIssues.Range("B6").Value = "Total number of unique key"
$wsIssues.Range("C6").Value = 0

$wsIssues.Name = "SCHEME_MEASURES"

# ---------------------------------------------------------------------------
# Sheet #3 : "2_SCHEME_MEASURES" -> "METADATA_ISSUES"
# Replace the measure counts with the (wider) issues table.
# ---------------------------------------------------------------------------
$wsScheme.Range("A2:C6").Clear()

$wsScheme.Range("A1").Value = "rule"
$wsScheme.Range("B1").Value = "desc"
$wsScheme.Range("C1").Value = "owner"
$wsScheme.Range("D1").Value = "table"
$wsScheme.Range("E1").Value = "column"
$wsScheme.Range("F1").Value = "constraint_name"
$wsScheme.Range("G1").Value = "length"
$wsScheme.Range("H1").Value = "limit"

$wsScheme.Range("A2").Value = "MQME01"
$wsScheme.Range("B2").Value = "Columns without comments"
$wsScheme.Range("C2").Value = "COPAF_AIX"
$wsScheme.Range("D2").Value = "DEBITO"
$wsScheme.Range("E2").Value = "COD_TIP_TRA_ITCD"
# The rule has no constraint/length/limit for this row; still realize the
# cells (empty, unstyled) instead of leaving them completely absent.
$wsScheme.Range("F2:H2").Borders.LineStyle = 0

$wsScheme.Name = "METADATA_ISSUES"

# ---------------------------------------------------------------------------
# Sheet #4 : "3_MODEL_MEASURES" -> "METADATA_MEASURES"
# Header (Indicator/Description/Value) stays the same; the single measure
# row is replaced and two more measure rows are appended.
# ---------------------------------------------------------------------------
$wsModelM.Range("A2").Value = "MQME00"
$wsModelM.Range("B2").Value = "Total number of columns"
$wsModelM.Range("C2").Value = 246

$wsModelM.Range("A3").Value = "MQMEA1"
$wsModelM.Range("B3").Value = "Total number of length-required columns"
$wsModelM.Range("C3").Value = 28

$wsModelM.Range("A4").Value = "MQMEA2"
$wsModelM.Range("B4").Value = "Total number of NUMBER columns"
$wsModelM.Range("C4").Value = 218

$wsModelM.Name = "METADATA_MEASURES"

# ---------------------------------------------------------------------------
# Sheet #5 : "4_MODEL_METRICS" -> "METADATA_METRICS"
# Header stays the same; indicator codes/labels/values are renumbered and
# the last (9th) row is dropped entirely. Values in column C are textual
# percentages ("99.59%"), so force text formatting before assigning them,
# then restore the plain "Normal" style so no stray number format lingers.
# ---------------------------------------------------------------------------
$wsModelMet.Range("A9:C9").Clear()

$wsModelMet.Range("A2").Value = "IQME01"
$wsModelMet.Range("B2").Value = "Columns with comments"
$wsModelMet.Range("C2").NumberFormat = "@"
$wsModelMet.Range("C2").Value = "99.59%"
$wsModelMet.Range("C2").Style = "Normal"

$wsModelMet.Range("A3").Value = "IQME02"
$wsModelMet.Range("B3").Value = "Columns with data type"
$wsModelMet.Range("C3").NumberFormat = "@"
$wsModelMet.Range("C3").Value = "100.00%"
$wsModelMet.Range("C3").Style = "Normal"

$wsModelMet.Range("A4").Value = "IQME03"
$wsModelMet.Range("B4").Value = "Length-required columns with data length"
$wsModelMet.Range("C4").NumberFormat = "@"
$wsModelMet.Range("C4").Value = "100.00%"
$wsModelMet.Range("C4").Style = "Normal"

$wsModelMet.Range("A5").Value = "IQME04"
$wsModelMet.Range("B5").Value = "NUMBER columns with valid scale"
$wsModelMet.Range("C5").NumberFormat = "@"
$wsModelMet.Range("C5").Value = "100.00%"
$wsModelMet.Range("C5").Style = "Normal"

$wsModelMet.Range("A6").Value = "IQME05"
$wsModelMet.Range("B6").Value = "Columns with valid num_distinct"
$wsModelMet.Range("C6").NumberFormat = "@"
$wsModelMet.Range("C6").Value = "100.00%"
$wsModelMet.Range("C6").Style = "Normal"

$wsModelMet.Range("A7").Value = "IQME06"
$wsModelMet.Range("B7").Value = "Columns with valid num_nulls"
$wsModelMet.Range("C7").NumberFormat = "@"
$wsModelMet.Range("C7").Value = "100.00%"
$wsModelMet.Range("C7").Style = "Normal"

$wsModelMet.Range("A8").Value = "IQME07"
$wsModelMet.Range("B8").Value = "Columns with valid density"
$wsModelMet.Range("C8").NumberFormat = "@"
$wsModelMet.Range("C8").Value = "100.00%"
$wsModelMet.Range("C8").Style = "Normal"

$wsModelMet.Name = "METADATA_METRICS"

Write-Host "done"
